# courses.xlsx — "1.07 with intent update"
#
# The "业委会工作" (HOA-committee-work) scenario row (row 9) is rewritten:
# its background description, character sketch and AI opening line are
# replaced with a new "committee re-election" scenario (徐阿姨 / 周老师),
# replacing the old "delivery locker" (快递柜 / 小苏) scenario. All other
# rows are untouched; any now-unused shared strings are dropped automatically
# when the workbook is saved.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("B9").Value = "某商品房小区正在进行业委会换届改选。业主徐阿姨提出小区建立二十多年来从未召开过业主大会，认为此次选举，居委会无视业主意见而完全包办，任命不符合要求的前业委会失职人员周老师继续做业委会候选人。徐阿姨对上述说辞其实并没有实际证据，之前因为拆违工作，她一直对居委会心怀不满，于是故意在选举投票现场公开对居委会发难。`n接下来我来扮演徐阿姨，而您扮演主持投票选举的工作人员。"

$ws.Range("C9").Value = "你是小区业主徐阿姨，今年68岁，国企退休职工，之前居委会强制拆除了你家的违章建筑，这让你很气愤，之后凡是居委会同意的你就反对。今天你们小区进行业委会换届改选投票，你注意到候选人中有前业委会成员周老师，你认定周老师之前有失职行为，且小区建立二十多年来从未召开过业主大会，你来到投票现场，大声向主持工作的居委工作人员发难，一定要让他下不来台！"

$ws.Range("D9").Value = "你们这次选举暗箱操作，谁给你们的权力选周老师，他不符合条件，你们是只手遮天。"

# Move the sheet's active cell / selection from C9 to C10, matching the
# author's final cursor position after editing the row above.
$null = $ws.Range("C10").Select()
